# fix resource url mapping for geojson and update dataset_json_lookup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("url" / "field_link_api" / FALSE) is a stray/incorrect mapping row.
# Delete the whole row so subsequent rows (country_code, region,
# license_title, id) shift up to close the gap.
$ws.Rows.Item(14).Delete()

# Mirror the author's final cursor position recorded in the saved file.
$ws.Range("D27").Select()
